$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Riscos")

# Update Probabilidade (column G) values; Magnitude (column H) holds the
# formula =F*G and recalculates automatically.
$ws.Range("G3").Value = 0.02
$ws.Range("G4").Value = 0.1
$ws.Range("G6").Value = 0.2
$ws.Range("G7").Value = 0.1

# Update the sheet's view state: scroll position and active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G7").Select()
